$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G13").Value  = 1.72
$ws.Range("I13").Value  = 4
$ws.Range("J13").Value  = 2.27
$ws.Range("K13").Value  = 2.25
$ws.Range("M13").Value  = 1.22
$ws.Range("N13").Value  = 3.4
$ws.Range("O13").Value  = 1.7
$ws.Range("P13").Value  = 1.93
$ws.Range("Q13").Value  = 2.6
$ws.Range("R13").Value  = 1.38
$ws.Range("U13").Value  = 1.65
$ws.Range("V13").Value  = 1.98
$ws.Range("W13").Value  = 7.9
$ws.Range("X13").Value  = 8.75
$ws.Range("AA13").Value = 13
$ws.Range("AC13").Value = 12.5
$ws.Range("AD13").Value = 7.5
$ws.Range("AF13").Value = 60
$ws.Range("AG13").Value = 450
$ws.Range("AH13").Value = 13
